$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 8379.375
$ws.Range("I74").Value = 5678.3335
$ws.Range("K74").Value = 5678.3335
$ws.Range("M74").Value = -4742.3335

# Row 77
$ws.Range("H77").Value = 8379.375
$ws.Range("I77").Value = 5678.3335
$ws.Range("K77").Value = 28391.6675
$ws.Range("M77").Value = -23711.6675

# Row 132
$ws.Range("H132").Value = 898.14703
$ws.Range("I132").Value = 512.8889
$ws.Range("J132").Value = 2384.1428
$ws.Range("K132").Value = 1538.6667
$ws.Range("L132").Value = 7152.428400000001
$ws.Range("M132").Value = 991.3332999999998
$ws.Range("N132").Value = -12212.4284

# Row 135
$ws.Range("H135").Value = 773.9394
$ws.Range("I135").Value = 616
$ws.Range("J135").Value = 1016.9231
$ws.Range("K135").Value = 5544
$ws.Range("L135").Value = 9152.3079
$ws.Range("M135").Value = -3009
$ws.Range("N135").Value = -14222.3079

# Row 137
$ws.Range("H137").Value = 6062.489
$ws.Range("I137").Value = 3846.75
$ws.Range("J137").Value = 11516.615
$ws.Range("K137").Value = 11540.25
$ws.Range("L137").Value = 34549.845
$ws.Range("M137").Value = -8990.25
$ws.Range("N137").Value = -39649.845

# Row 138
$ws.Range("H138").Value = 4063.4792
$ws.Range("I138").Value = 3911.175
$ws.Range("J138").Value = 4825
$ws.Range("K138").Value = 11733.525
$ws.Range("L138").Value = 14475
$ws.Range("M138").Value = -6593.525000000001
$ws.Range("N138").Value = -24755

# Row 141
$ws.Range("H141").Value = 2263.0557
$ws.Range("I141").Value = 2140
$ws.Range("J141").Value = 3247.5
$ws.Range("K141").Value = 6420
$ws.Range("L141").Value = 9742.5
$ws.Range("M141").Value = -1240
$ws.Range("N141").Value = -20102.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16530.1
$ws.Range("I32").Value = 15927.39
$ws.Range("J32").Value = 31999.666
$ws.Range("K32").Value = 15927.39
$ws.Range("L32").Value = 31999.666
$ws.Range("M32").Value = -15640.39
$ws.Range("N32").Value = -32573.666

# Row 74
$ws.Range("H74").Value = 347765.75
$ws.Range("I74").Value = 502187.06
$ws.Range("J74").Value = 4607.3335
$ws.Range("K74").Value = 502187.06
$ws.Range("L74").Value = 4607.3335
$ws.Range("M74").Value = -501313.06
$ws.Range("N74").Value = -6355.3335

# Row 77
$ws.Range("H77").Value = 347765.75
$ws.Range("I77").Value = 502187.06
$ws.Range("J77").Value = 4607.3335
$ws.Range("K77").Value = 2510935.3
$ws.Range("L77").Value = 23036.6675
$ws.Range("M77").Value = -2506567.3
$ws.Range("N77").Value = -31772.6675

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 946.96295
$ws.Range("I107").Value = 933.73914
$ws.Range("J107").Value = 1023
$ws.Range("K107").Value = 933.73914
$ws.Range("L107").Value = 1023
$ws.Range("M107").Value = 986.26086
$ws.Range("N107").Value = -4863

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4570.891
$ws.Range("I31").Value = 1900.3636
$ws.Range("J31").Value = 5410.2
$ws.Range("K31").Value = 1900.3636
$ws.Range("L31").Value = 5410.2
$ws.Range("M31").Value = -1605.3636
$ws.Range("N31").Value = -6000.2

# Row 34
$ws.Range("H34").Value = 4570.891
$ws.Range("I34").Value = 1900.3636
$ws.Range("J34").Value = 5410.2
$ws.Range("K34").Value = 1900.3636
$ws.Range("L34").Value = 5410.2
$ws.Range("M34").Value = -1698.3636
$ws.Range("N34").Value = -5814.2

# Row 99
$ws.Range("H99").Value = 7399.6
$ws.Range("I99").Value = 8332.666999999999
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 8332.666999999999
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -6834.666999999999
$ws.Range("N99").Value = -8996

# Row 126
$ws.Range("H126").Value = 7399.6
$ws.Range("I126").Value = 8332.666999999999
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 24998.001
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -22528.001
$ws.Range("N126").Value = -22940

# Row 132
$ws.Range("H132").Value = 3596.75
$ws.Range("I132").Value = 2675.4092
$ws.Range("J132").Value = 6975
$ws.Range("K132").Value = 8026.2276
$ws.Range("L132").Value = 20925
$ws.Range("M132").Value = -5496.2276
$ws.Range("N132").Value = -25985

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = $null
$ws.Range("N134").Value = $null

# Row 140
$ws.Range("H140").Value = 119999
$ws.Range("J140").Value = 119999
$ws.Range("L140").Value = 119999
$ws.Range("N140").Value = -130359

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 2303
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2303
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 9212
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -9350

# Row 58
$ws.Range("H58").Value = 2702
$ws.Range("I58").Value = 3005
$ws.Range("K58").Value = 9015
$ws.Range("M58").Value = -8887

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5134.5625
$ws.Range("I80").Value = 4198.6
$ws.Range("J80").Value = 5560
$ws.Range("K80").Value = 4198.6
$ws.Range("L80").Value = 5560
$ws.Range("M80").Value = -3200.6
$ws.Range("N80").Value = -7556

# Row 83
$ws.Range("H83").Value = 5134.5625
$ws.Range("I83").Value = 4198.6
$ws.Range("J83").Value = 5560
$ws.Range("K83").Value = 20993
$ws.Range("L83").Value = 27800
$ws.Range("M83").Value = -16001
$ws.Range("N83").Value = -37784

# Row 122
$ws.Range("H122").Value = 2012.0869
$ws.Range("I122").Value = 1721.0667
$ws.Range("J122").Value = 2557.75
$ws.Range("K122").Value = 5163.2001
$ws.Range("L122").Value = 7673.25
$ws.Range("M122").Value = -2713.2001
$ws.Range("N122").Value = -12573.25

# Row 132
$ws.Range("H132").Value = 3037.352
$ws.Range("I132").Value = 2399
$ws.Range("J132").Value = 13729.75
$ws.Range("K132").Value = 7197
$ws.Range("L132").Value = 41189.25
$ws.Range("M132").Value = -4667
$ws.Range("N132").Value = -46249.25

$ws = $wb.Worksheets.Item("LTW")
# Row 137
$ws.Range("H137").Value = 126666.5
$ws.Range("J137").Value = 126666.5
$ws.Range("L137").Value = 126666.5
$ws.Range("N137").Value = -136866.5

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1788.9688
$ws.Range("I122").Value = 1770.7142
$ws.Range("J122").Value = 1916.75
$ws.Range("K122").Value = 5312.142599999999
$ws.Range("L122").Value = 5750.25
$ws.Range("M122").Value = -2862.142599999999
$ws.Range("N122").Value = -10650.25

# Row 124
$ws.Range("H124").Value = 113143
$ws.Range("J124").Value = 113143
$ws.Range("L124").Value = 113143
$ws.Range("N124").Value = -122963

# Row 126
$ws.Range("H126").Value = 7415.3076
$ws.Range("I126").Value = 7562.143
$ws.Range("J126").Value = 7244
$ws.Range("K126").Value = 22686.429
$ws.Range("L126").Value = 21732
$ws.Range("M126").Value = -20216.429
$ws.Range("N126").Value = -26672
